$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel
# need an explicit text format so they stay text (matching the source data,
# which stores every Price/Volume cell as a string).

$ws.Range("D2").Value = '26.073.36'
$ws.Range("E2").Value = '  +5.74%  '
$ws.Range("D3").Value = '1.717.07'
$ws.Range("E3").Value = '  +3.60%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '333.99'
$ws.Range("E5").Value = '  +4.16%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.13%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3688'
$ws.Range("E7").Value = '  +1.30%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '49.30'
$ws.Range("E8").Value = '  +5.23%  '
$ws.Range("E10").Value = '  +4.53%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07465'
$ws.Range("E11").Value = '  +5.87%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.002'
$ws.Range("E12").Value = '  +0.16%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.295'
$ws.Range("E13").Value = '  +5.21%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.931'
$ws.Range("E15").Value = '  +4.71%  '
$ws.Range("D16").Value = '1.720.32'
$ws.Range("E16").Value = '  +3.75%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001076'
$ws.Range("E17").Value = '  +2.89%  '
$ws.Range("E18").Value = '  +0.15%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '81.94'
$ws.Range("E19").Value = '  +3.89%  '
$ws.Range("E20").Value = '  +0.17%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.39'
$ws.Range("E21").Value = '  +4.12%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.088'
$ws.Range("E22").Value = '  +2.81%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.02'
$ws.Range("E23").Value = '  +3.20%  '
$ws.Range("D24").Value = '26.022.43'
$ws.Range("E24").Value = '  +5.58%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.475'
$ws.Range("E25").Value = '  +0.94%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.453'
$ws.Range("E26").Value = '  +2.67%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '150.52'
$ws.Range("E27").Value = '  +1.71%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.21'
$ws.Range("E28").Value = '  +3.30%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.318'
$ws.Range("E29").Value = '  +8.95%  '
$ws.Range("D30").Value = '1.909.08'
$ws.Range("E30").Value = '  +3.65%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '129.31'
$ws.Range("E31").Value = '  +3.03%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.124'
$ws.Range("E32").Value = '  +1.23%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.917'
$ws.Range("E33").Value = '  +1.03%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '12.87'
$ws.Range("E36").Value = '  +4.59%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.347'
$ws.Range("E37").Value = '  +2.50%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06226'
$ws.Range("E38").Value = '  +3.27%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02292'
$ws.Range("E39").Value = '  +2.55%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2134'
$ws.Range("E40").Value = '  +2.77%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.548'
$ws.Range("E41").Value = '  +4.37%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.221'
$ws.Range("E42").Value = '  -4.55%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '14.51'
$ws.Range("E43").Value = '  +13.28%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6166'
$ws.Range("E44").Value = '  +4.05%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.001'
$ws.Range("E45").Value = '  +0.20%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.834'
$ws.Range("E46").Value = '  -0.48%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5890'
$ws.Range("E47").Value = '  +4.65%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '128.04'
$ws.Range("E48").Value = '  +2.84%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.019'
$ws.Range("E49").Value = '  +3.34%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07263'
$ws.Range("E50").Value = '  +4.28%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '77.05'
$ws.Range("E51").Value = '  +3.25%  '

# Row 34 <-> 35 swap (Stellar and WEMIXTOKEN exchange rank positions)
$ws.Range("D34").NumberFormat = "@"
$ws.Range("B34").Value = 'WEMIXTOKEN'
$ws.Range("C34").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D34").Value = '1.722'
$ws.Range("E34").Value = '  +2.75%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("B35").Value = 'Stellar'
$ws.Range("C35").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D35").Value = '0.08512'
$ws.Range("E35").Value = '  +0.48%  '
